$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column N (27-jun) ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header for new column, matching style of existing header cells (e.g. M1)
$wsPrix.Range("M1").Copy()
$wsPrix.Range("N1").PasteSpecial(-4122)  # xlPasteFormats
$wsPrix.Range("N1").Value = "27-jun"

$nValues = @{
    2  = 81.34
    3  = 66.17
    4  = 63.79
    5  = 45.14
    6  = 46.62
    7  = 71.3
    8  = 84.53
    9  = 90.09
    10 = 73.55
    11 = 37.04
    12 = 1.34
    13 = 0
    14 = -0.01
    15 = -0.01
    16 = -0.03
    17 = -0.01
    18 = 7.5
    19 = 20.06
    20 = 71.95
    21 = 104.61
    22 = 125.4
    23 = 125.49
    24 = 131.91
    25 = 102.48
}

foreach ($row in $nValues.Keys) {
    $wsPrix.Cells.Item($row, 14).Value = $nValues[$row]
}

# --- Sheet "Gaz": update B8 value ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("B8").Value = 34.75
